# Regenerate the "K" column (column G) values for save_data.
# The commit replaces the old "Strike#" derived values with the new "K" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value (column G), per the recomputed s_vals.
$kValues = @{
    2  = 8
    3  = 9
    4  = 7
    5  = 4
    6  = 7
    7  = 6
    8  = 9
    9  = 6
    10 = 8
    11 = 7
    12 = 10
    13 = 6
    14 = 6
    15 = 6
    16 = 8
    17 = 4
    18 = 9
    19 = 5
    20 = 10
    21 = 8
    22 = 6
    23 = 9
    24 = 6
    25 = 6
    26 = 9
    27 = 12
    28 = 6
    29 = 5
    31 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
